# Apply "Add data for 2022-09-06" update to the carjacking YoY workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet and update its on-sheet label reflecting the new "through" date.
$ws.Name = "Through 2022-08-29"

# Update the shared string label for the August row (column A, row 9).
$ws.Range("A9").Value = "August (through 08-29)"

# Update August row (row 9) values.
$ws.Range("B9").Value = 30
$ws.Range("C9").Value = 74
$ws.Range("D9").Value = 84
$ws.Range("E9").Value = 62
$ws.Range("F9").Value = 42
$ws.Range("G9").Value = 155
$ws.Range("H9").Value = 152
$ws.Range("I9").Value = 152

# Update Total row (row 10) values.
$ws.Range("B10").Value = 192
$ws.Range("C10").Value = 376
$ws.Range("D10").Value = 549
$ws.Range("E10").Value = 487
$ws.Range("F10").Value = 346
$ws.Range("G10").Value = 776
$ws.Range("H10").Value = 1062
$ws.Range("I10").Value = 1123
